# The workbook's monthly data is grouped in 12-row blocks, one per year
# (rows 2-13 = 2014, 14-25 = 2015, 26-37 = 2016, 38-49 = 2017). The edit
# reorders each year block so that the Oct/Nov/Dec rows move to the front
# of the block (ahead of Jan), i.e. the last 3 rows of every 12-row block
# are rotated to the top: Oct,Nov,Dec,Jan,Feb,...,Sep.
#
# We implement this generically: for every 12-row block, read all the
# existing rows into memory first (so we never read data that has
# already been overwritten), then write them back out in rotated order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 49
$blockSize = 12
$numCols = 5

for ($blockStart = $firstDataRow; $blockStart -le $lastDataRow; $blockStart += $blockSize) {
    $blockEnd = $blockStart + $blockSize - 1

    # Read every row of this block into memory before writing anything back.
    $rowsData = @{}
    for ($r = $blockStart; $r -le $blockEnd; $r++) {
        $vals = @()
        for ($c = 1; $c -le $numCols; $c++) {
            $vals += , ($ws.Cells.Item($r, $c).Value2)
        }
        $rowsData[$r] = $vals
    }

    # Rotate: the last 3 rows of the block move to the front, the first
    # 9 rows shift down by 3 positions.
    $rotated = @()
    for ($i = 9; $i -le 11; $i++) {
        $rotated += , ($blockStart + $i)
    }
    for ($i = 0; $i -le 8; $i++) {
        $rotated += , ($blockStart + $i)
    }

    $destRow = $blockStart
    foreach ($srcRow in $rotated) {
        $vals = $rowsData[$srcRow]
        for ($c = 1; $c -le $numCols; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
        }
        $destRow++
    }
}
